$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.588.47'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '1.665.32'
$ws.Range('E3').Value = '  -3.56%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.17'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.53'
$ws.Range('E8').Value = '  -2.71%  '
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0621'
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('E12').Value = '  -3.51%  '
$ws.Range('D13').Value = '1.687.63'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  -2.66%  '
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.17'
$ws.Range('E16').Value = '  -2.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '248.24'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').Value = '27.607.35'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('E20').Value = '  -4.55%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.49'
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.29'
$ws.Range('E23').Value = '  -4.98%  '
$ws.Range('E24').Value = '  -5.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.09'
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.16'
$ws.Range('E26').Value = '  -4.97%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.40'
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('E30').Value = '  +3.58%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('D33').Value = '1.474.12'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('E35').Value = '  -5.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.936'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.574'
$ws.Range('E38').Value = '  -6.21%  '
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.58'
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('E41').Value = '  -5.45%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.39'
$ws.Range('E43').Value = '  -7.68%  '
$ws.Range('D44').Value = '1.809.62'
$ws.Range('E44').Value = '  -3.47%  '
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.788'
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.70'
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.31'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('D49').Value = '0.0₆0109'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.90'
$ws.Range('E51').Value = '  -4.00%  '
